$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G21").Value = "RS"
